$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 37785982911.24337
$ws.Range("H2").Value = 1675881703641.516
$ws.Range("I2").Value = 28754400962834900

# Row 3
$ws.Range("B3").Value = 3.045061224774464
$ws.Range("C3").Value = 496.1349728514716
$ws.Range("D3").Value = 51463.44232116333
$ws.Range("E3").Value = 3949555.41930737
$ws.Range("F3").Value = 302029135.8022056
$ws.Range("G3").Value = 12023767741.14577
$ws.Range("H3").Value = 592496493073.3569
$ws.Range("I3").Value = 12732581342679240

# Row 4
$ws.Range("B4").Value = 8.663526526158217
$ws.Range("C4").Value = 1843.119817584459
$ws.Range("D4").Value = 239029.2143902132
$ws.Range("E4").Value = 16402802.7359937
$ws.Range("F4").Value = 1100370138.151046
$ws.Range("G4").Value = 38817825716.94899
$ws.Range("H4").Value = 1712728393786.48
$ws.Range("I4").Value = 29150968857903340

# Row 5
$ws.Range("G5").Value = 36559212319.26645
$ws.Range("H5").Value = 1620746478918.832
$ws.Range("I5").Value = 28702183129570312

# Row 6
$ws.Range("B6").Value = 5.369242002566493
$ws.Range("C6").Value = 887.1610385335186
$ws.Range("D6").Value = 76725.4774647124
$ws.Range("E6").Value = 4876168.263484871
$ws.Range("F6").Value = 373799034.0554756
$ws.Range("G6").Value = 13455285178.20615
$ws.Range("H6").Value = 654006570286.0819
$ws.Range("I6").Value = 10839406375384220

# Row 7
$ws.Range("G7").Value = 1287459788.280791
$ws.Range("H7").Value = 64380602828.17851
$ws.Range("I7").Value = 1476931647182352

# Row 8
$ws.Range("G8").Value = 8148610972.51496
$ws.Range("H8").Value = 405391373034.7842
$ws.Range("I8").Value = 8995373167014003

# Row 9
$ws.Range("B9").Value = 2.701955361016653
$ws.Range("C9").Value = 607.6788956355718
$ws.Range("D9").Value = 94891.83174998802
$ws.Range("E9").Value = 7650016.09706249
$ws.Range("F9").Value = 574225765.8906945
$ws.Range("G9").Value = 22120957776.42364
$ws.Range("H9").Value = 1052859501170.646
$ws.Range("I9").Value = 20794986822359088

# Row 10
$ws.Range("G10").Value = 1400.979036837564
$ws.Range("H10").Value = 97372.7654354666
$ws.Range("I10").Value = 1502699451.718121
